$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in missing attendance marks for several students
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 5
$ws.Range("C9").Value = 5
$ws.Range("D9").Value = 5
$ws.Range("D16").Value = 5
$ws.Range("D27").Value = 5

# Move the active selection to D29 (also updates the frozen-pane top-left cell)
$ws.Range("D29").Select()
